$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B1").Value = 1
$ws.Range("B2").Value = 2
$ws.Range("B3").Value = 3
$ws.Range("B4").Value = 4

$ws.Range("B4").Select()
